$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking text (e.g. "1.000", "0.1000") that must
# stay as literal text instead of being auto-coerced into numbers, so format the
# whole price column as Text before writing any values into it.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.564.91'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.875.70'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '247.15'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').Value = '0.9995'
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').Value = '0.4741'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D8').Value = '0.2918'
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('D9').Value = '0.06487'
$ws.Range('E9').Value = '  +0.26%  '
$ws.Range('D10').Value = '22.13'
$ws.Range('E10').Value = '  +4.99%  '
$ws.Range('D11').Value = '0.07722'
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '0.7466'
$ws.Range('E12').Value = '  +5.41%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').Value = '97.13'
$ws.Range('E13').Value = '  +2.18%  '
$ws.Range('D14').Value = '1.872.15'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '5.156'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '274.24'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '30.554.00'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('D18').Value = '13.42'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000007521'
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').Value = '2.116.84'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('D23').Value = '5.282'
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').Value = '6.185'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').Value = '9.269'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').Value = '163.64'
$ws.Range('E26').Value = '  -0.65%  '
$ws.Range('D27').Value = '18.83'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = '1.921'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('D29').Value = '0.1000'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('E30').Value = '  -1.74%  '
$ws.Range('D31').Value = '1.511'
$ws.Range('E31').Value = '  -0.48%  '
$ws.Range('D32').Value = '4.301'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').Value = '4.126'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('D34').Value = '0.04810'
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').Value = '0.6989'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '0.9997'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('B38').Value = 'HuobiToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01846'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.750'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '6.234'
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '74.11'
$ws.Range('E42').Value = '  +5.13%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.973'
$ws.Range('E43').Value = '  +3.64%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4184'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9999'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '0.8346'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = '102.78'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.385'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').Value = '937.89'
$ws.Range('E49').Value = '  +2.06%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '35.37'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').Value = '6.979'
$ws.Range('E51').Value = '  -1.40%  '
